$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "x" flag and bill amount from row 19 up to row 16 (AIA forms),
# casting the bill amount to a float.
$ws.Range("B16").Value = "x"
$ws.Range("D16").Value = [double]340.04

$ws.Range("B19").ClearContents()
$ws.Range("D19").ClearContents()

# Update the running invoice/date counter in I1.
$ws.Range("I1").Value = 33453

# Reset the active selection to I1.
$ws.Range("I1").Select()
